$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "<his>"

$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 11
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 4
$ws.Range("C9").Value = 6
$ws.Range("C10").Value = 11
$ws.Range("C13").Value = 10
$ws.Range("C16").Value = 4
